# Apply the recorded edit to the first (only) paragraph of the document:
#   1. Give the paragraph the "Abstract" paragraph style.
#   2. Duplicate its single " " run so the paragraph holds two runs of " ".
#   3. Relocate the "_GoBack" bookmark so it wraps an empty range placed
#      after both runs (i.e. at the very end of the paragraph's text),
#      instead of wrapping the whole paragraph content as before.

$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)

# 1. Set the paragraph style to "Abstract".
$p.Style = "Abstract"

# 2. Append a duplicate run containing a single space, matching the
#    existing run's text, right after the existing content.
$r = $p.Range
$r.InsertAfter(" ")

# Temporarily append a sentinel character after the two runs so that the
# insertion point used for the bookmark is NOT touching the paragraph
# mark; inserting/collapsing a bookmark range right at the paragraph's
# trailing boundary causes it to wrap the entire paragraph instead of
# sitting after the runs, so we avoid that boundary case.
$r2 = $p.Range
$r2.InsertAfter("X")

# Remove the old bookmark (currently still wrapping the original run at
# the start of the paragraph) so it can be recreated in its new spot.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 3. Recreate the "_GoBack" bookmark as a zero-length range positioned
#    right after the two space runs (and before the sentinel character).
$paraStart = $p.Range.Start
$bmRange = $d.Range($paraStart + 2, $paraStart + 2)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the sentinel character now that the bookmark is anchored
# correctly; the bookmark remains after both runs.
$delRange = $d.Range($paraStart + 2, $paraStart + 3)
$delRange.Delete()
